$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Yashshree", 2, "chinmaymhatre21@gmail.com", "19/11/2024, 6:41:05 pm", "34734b0c"),
    @("Khushil", 2, "khushilbhimani2@gmail.com", "19/11/2024, 6:53:00 pm", "c881e563"),
    @("Khushil", 2, "khushilbhimani2@gmail.com", "19/11/2024, 7:01:01 pm", "80a40d9e"),
    @("Khushil", 2, "khushilbhimani2@gmail.com", "19/11/2024, 7:04:46 pm", "122441d0"),
    @("Khushil Girish Bhimani", 2, "khushilbhimani302@gmail.com", "19/11/2024, 7:06:07 pm", "ecb8c446"),
    @("Khushil Bhimani", 2, "khushilbhimani2@gmail.com", "19/11/2024, 7:13:09 pm", "555e4dae"),
    @("Khush", 2, "khushilbhimani2@gmail.com", "19/11/2024, 7:14:38 pm", "19c59d21"),
    @("Khushil", 2, "Khushilbhimani2@gmail.com", "19/11/2024, 7:16:28 pm", "c9e25a3c"),
    @("Khushil Bhimani", 2, "khushilbhimani302@gmail.com", "19/11/2024, 7:19:07 pm", "6cd1327d"),
    @("Khuhsil", 2, "Khushilbhimani302@gmail.com", "19/11/2024, 7:20:14 pm", "7ab50f23"),
    @("Khushil Girish Bhimani", 2, "khushilbhimani302@gmail.com", "19/11/2024, 7:21:18 pm", "4e6b9388"),
    @("pradipgarhwal183@gmail.com", 2, "pradipgarhwal183@gmail.com", "19/11/2024, 7:40:59 pm", "0dbacf79"),
    @("Chinmay Mhatre", 2, "chinmaymhatre21@gmail.com", "21/11/2024, 12:28:26 am", "9878e1fd"),
    @("Pradip Garhwal", 2, "pradipgarhwal183@gmail.com", "21/11/2024, 12:30:20 am", "99e18d11"),
    @("Rohit Gautam", 0, "rohit@hacktify.in", "21/11/2024, 12:44:09 am", "2a846fa0")
)

$row = 9
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
